$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFile = "df4e5c99-e057-4b1b-bb78-837172eff53c.md"
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37395a8d858c0ddaf12e43f43f4e14dede6a6c0f/e2e/df4e5c99-e057-4b1b-bb78-837172eff53c.md"

# --- Overview sheet: update Status for both languages ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# --- zh-cn sheet: record the handback results ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("I2").Value = $targetFile
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $handbackUrl, [Type]::Missing, [Type]::Missing, $targetFile) | Out-Null
$wsZh.Range("J2").Value = "df4e5c99-e057-4b1b-bb78-837172eff53c.b6451b5773a334ee600a3d527af2cc6c67f8512f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 01:03:28"

# --- de-de sheet: record the handback results ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("I2").Value = $targetFile
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $handbackUrl, [Type]::Missing, [Type]::Missing, $targetFile) | Out-Null
$wsDe.Range("J2").Value = "df4e5c99-e057-4b1b-bb78-837172eff53c.b6451b5773a334ee600a3d527af2cc6c67f8512f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 01:03:35"

# --- Column width refresh (best-effort autofit for widened status/report columns) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17
$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
